$p = $ppt.ActivePresentation

# 1. Add a new paragraph to the "Implementation" slide (slide 4) content placeholder.
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.InsertAfter("`rgit was used for version control")

# 2. Add a new "Conclusion" slide at the end, using the Title and Content layout.
$newSlide = $p.Slides.Add(5, 2)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusion"
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Project functions as initially desired`rLearned a lot about socket programming and multi-threading`rThanks for listening"
